$wb = $excel.ActiveWorkbook

# --- "Data" sheet: the browser value in C2 changes from "chromegrid" to "chrome" ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("C2").Value = "chrome"

# --- "Test" sheet: selection becomes a multi-area pick of C2 and D8, with D8 active ---
$wsTest = $wb.Worksheets.Item("Test")
[void]$wsTest.Activate()
[void]$wsTest.Range("C2,D8").Select()
[void]$wsTest.Range("D8").Activate()

# --- "Data" sheet: becomes the active sheet again, with C2 selected ---
[void]$wsData.Activate()
[void]$wsData.Range("C2").Select()
